# Updated Argent prices in Excel
# Appends a new "2025-03-04" price row (row 3) to the sheets that only
# had data through 2025-03-03 (row 2), carrying forward each sheet's
# latest price value. Cell values are written as plain text (matching
# the workbook's existing inlineStr / text-typed cells), not as Excel
# dates or numbers.

$wb = $excel.ActiveWorkbook

function Add-PriceRow {
    param(
        [string]$SheetName,
        [string]$Date,
        [string]$Price
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Force the new values to be stored as text (leading apostrophe),
    # then reset the style so no number-format/quote-prefix style is
    # left attached to the cells.
    $ws.Range("A3").Formula = "'" + $Date
    $ws.Range("B3").Formula = "'" + $Price
    $ws.Range("A3:B3").Style = "Normal"
}

Add-PriceRow "N-Dense" "2025-03-04" "5.48"
Add-PriceRow "N-Type"  "2025-03-04" "5.89"
Add-PriceRow "USD_CNY" "2025-03-04" "7.3048"
